# Add data for 2022-11-28
# Update year-to-date violent crime counts across the citywide totals,
# the "By Neighborhood" rollup, and the affected individual neighborhood
# sheets to reflect the new day of data.

$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("D3").Value = 131
$ws.Range("I3").Value = 188
$ws.Range("C6").Value = 466
$ws.Range("D6").Value = 401
$ws.Range("E6").Value = 458
$ws.Range("F6").Value = 510
$ws.Range("H6").Value = 434
$ws.Range("I6").Value = 490
$ws.Range("C7").Value = 620
$ws.Range("D7").Value = 629
$ws.Range("E7").Value = 680
$ws.Range("F7").Value = 739
$ws.Range("H7").Value = 703
$ws.Range("I7").Value = 817

# Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("C6").Value = 32
$ws.Range("C7").Value = 37

# By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("D7").Value = 10
$ws.Range("F8").Value = 49
$ws.Range("I8").Value = 40
$ws.Range("E19").Value = 15
$ws.Range("C36").Value = 37
$ws.Range("F53").Value = 76
$ws.Range("I53").Value = 122
$ws.Range("E54").Value = 5
$ws.Range("D76").Value = 14
$ws.Range("F77").Value = 20
$ws.Range("H77").Value = 29
$ws.Range("F85").Value = 5
$ws.Range("D96").Value = 10
$ws.Range("C98").Value = 620
$ws.Range("D98").Value = 629
$ws.Range("E98").Value = 680
$ws.Range("F98").Value = 739
$ws.Range("H98").Value = 703
$ws.Range("I98").Value = 817

# Woodlawn
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("D5").Value = 5
$ws.Range("D6").Value = 10

# Loop
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 29
$ws.Range("F6").Value = 56
$ws.Range("F7").Value = 76
$ws.Range("I7").Value = 122

# Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("D3").Value = 3
$ws.Range("D6").Value = 14

# Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("F6").Value = 10
$ws.Range("H6").Value = 18
$ws.Range("F7").Value = 20
$ws.Range("H7").Value = 29

# Chatham
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("E5").Value = 12
$ws.Range("E6").Value = 15

# Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("D5").Value = 7
$ws.Range("D6").Value = 10

# United Center
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = 5

# Lower West Side
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("E4").Value = 5
$ws.Range("E5").Value = 5

# Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("F5").Value = 35
$ws.Range("I5").Value = 30
$ws.Range("F6").Value = 49
$ws.Range("I6").Value = 40
